# Slide 4 ("Exam 1"), body placeholder, 2nd paragraph:
#   "1. \u201cWritten\u201d part on Moodle.  Closed book, notes, etc.  25-40%"
# becomes three runs:
#   "1" | ". Written " | "part.  Closed book, notes, etc.  25-40%"
# i.e. the smart-quotes are dropped and " on Moodle" is removed, and the
# run is split into three pieces around that edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(2)

# Sanity check we found the right paragraph before editing it.
if ($para.Text.IndexOf("Written") -ge 0) {

    # Drop the curly quotes around Written -> Written
    $quoted = $para.Characters(4, 9)
    $quoted.Text = "Written"

    # Drop " on Moodle"
    $cur = $para.Text
    $idx = $cur.IndexOf(" on Moodle")
    if ($idx -ge 0) {
        $onMoodle = $para.Characters($idx + 1, 10)
        $onMoodle.Text = ""
    }

    # Re-assert the three final runs so the paragraph ends up split as
    # "1" / ". Written " / "part.  Closed book, notes, etc.  25-40%"
    $run1 = $para.Characters(1, 1)
    $run1.Text = "1"

    $run2 = $para.Characters(2, 10)
    $run2.Text = ". Written "

    $run3 = $para.Characters(12, 39)
    $run3.Text = "part.  Closed book, notes, etc.  25-40%"
}
